$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to remain text (matches source inlineStr type,
# and avoids Excel auto-converting numeric-looking strings to floating point
# numbers which would lose trailing zeros / introduce FP noise).
$priceCells = @("D2","D3","D4","D6","D7","D8","D9","D10","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D28","D29","D30","D31","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.458.56"
$ws.Range("E2").Value = "  -0.68%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.832.62"
$ws.Range("E3").Value = "  +1.09%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.9966"
$ws.Range("E4").Value = "  -0.68%  "

# Row 5 - BNB
$ws.Range("E5").Value = "  -0.61%  "

# Row 6 - USDC
$ws.Range("D6").Value = "0.9934"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4469"
$ws.Range("E7").Value = "  +1.83%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3792"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9 - OKB
$ws.Range("D9").Value = "45.43"
$ws.Range("E9").Value = "  +1.96%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.07814"
$ws.Range("E10").Value = "  +1.05%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  -0.68%  "

# Row 12 - Solana
$ws.Range("E12").Value = "  -2.06%  "

# Row 13 - BinanceUSD
$ws.Range("D13").Value = "0.9954"
$ws.Range("E13").Value = "  -0.55%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.347"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "7.562"
$ws.Range("E15").Value = "  -0.79%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.830.21"
$ws.Range("E16").Value = "  +0.95%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "92.50"
$ws.Range("E17").Value = "  +13.48%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.00001088"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -5.33%  "

# Row 20 - Dai
$ws.Range("D20").Value = "0.9930"
$ws.Range("E20").Value = "  -0.66%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "17.61"
$ws.Range("E21").Value = "  -0.73%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "6.394"
$ws.Range("E22").Value = "  +1.12%  "

# Row 23 - BitDAO
$ws.Range("D23").Value = "0.5391"
$ws.Range("E23").Value = "  -0.82%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "28.507.46"
$ws.Range("E24").Value = "  -0.49%  "

# Row 25 - Cosmos
$ws.Range("D25").Value = "11.86"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "2.211"
$ws.Range("E26").Value = "  -9.38%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  +1.52%  "

# Row 28 - Monero
$ws.Range("D28").Value = "154.34"
$ws.Range("E28").Value = "  +0.64%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").Value = "2.386"
$ws.Range("E29").Value = "  +0.21%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "2.036.85"
$ws.Range("E30").Value = "  +0.84%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "130.19"
$ws.Range("E31").Value = "  -1.94%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "1.226"
$ws.Range("E32").Value = "  -3.94%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "5.884"
$ws.Range("E33").Value = "  +0.47%  "

# Row 34 - Stellar
$ws.Range("D34").Value = "0.09303"
$ws.Range("E34").Value = "  +0.13%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "3.676"
$ws.Range("E35").Value = "  -7.53%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "12.94"
$ws.Range("E36").Value = "  +5.38%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +0.73%  "

# Row 38 - Algorand
$ws.Range("D38").Value = "0.2204"
$ws.Range("E38").Value = "  -2.92%  "

# Row 39 - TheSandbox
$ws.Range("D39").Value = "0.6667"
$ws.Range("E39").Value = "  -0.17%  "

# Row 40 - InternetComputer(DFINITY)
$ws.Range("D40").Value = "5.215"
$ws.Range("E40").Value = "  -0.37%  "

# Row 41 - Hedera
$ws.Range("D41").Value = "0.06280"
$ws.Range("E41").Value = "  -1.64%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value = "1.195"
$ws.Range("E42").Value = "  -1.31%  "

# Row 43 - FraxShare
$ws.Range("D43").Value = "8.116"
$ws.Range("E43").Value = "  -0.68%  "

# Row 44 - WEMIXTOKEN
$ws.Range("D44").Value = "1.403"
$ws.Range("E44").Value = "  -3.39%  "

# Row 45 & 46 - Frax and EnergySwap swap ranking positions
# Row 45 becomes Frax (was EnergySwap), Row 46 becomes EnergySwap (was Frax)
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "0.9935"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "14.01"
$ws.Range("E46").Value = "  -0.23%  "

# Row 47 - Decentraland
$ws.Range("D47").Value = "0.6137"
$ws.Range("E47").Value = "  +0.28%  "

# Row 48 - PancakeSwap
$ws.Range("D48").Value = "3.763"
$ws.Range("E48").Value = "  -1.37%  "

# Row 49 - Quant
$ws.Range("D49").Value = "127.65"
$ws.Range("E49").Value = "  -1.11%  "

# Row 50 - NEARProtocol
$ws.Range("D50").Value = "2.048"
$ws.Range("E50").Value = "  -0.31%  "

# Row 51 - Aave
$ws.Range("D51").Value = "79.85"
$ws.Range("E51").Value = "  +1.58%  "
